$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("enc_mpr")

$ws.Range("A2").Value = "e"
$ws.Range("B2").Value = 181

$ws.Range("A3").Value = "."
$ws.Range("B3").Value = 32

$ws.Range("A4").Value = "g"
$ws.Range("B4").Value = 25

$ws.Range("A5").Value = "`n"
$ws.Range("B5").Value = 3

$ws.Range("A6").Value = "M"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "n"
$ws.Range("B7").Value = 82

$ws.Range("A8").Value = "l"
$ws.Range("B8").Value = 75

$ws.Range("A9").Value = "D"
$ws.Range("B9").Value = 3

$ws.Range("A10").Value = "s"
$ws.Range("B10").Value = 143

$ws.Range("A11").Value = "u"
$ws.Range("B11").Value = 122

$ws.Range("A12").Value = "t"
$ws.Range("B12").Value = 130

$ws.Range("A13").Value = "q"
$ws.Range("B13").Value = 9

$ws.Range("A14").Value = "m"
$ws.Range("B14").Value = 47

$ws.Range("A15").Value = "c"
$ws.Range("B15").Value = 61

$ws.Range("A16").Value = "a"
$ws.Range("B16").Value = 79

$ws.Range("A17").Value = ","
$ws.Range("B17").Value = 19

$ws.Range("A18").Value = "N"
$ws.Range("B18").Value = 4

$ws.Range("A19").Value = "F"
$ws.Range("B19").Value = 3

$ws.Range("A20").Value = "P"
$ws.Range("B20").Value = 2

$ws.Range("A21").Value = "d"
$ws.Range("B21").Value = 49

$ws.Range("A22").Value = "U"
$ws.Range("B22").Value = 1

$ws.Range("A23").Value = "h"
$ws.Range("B23").Value = 13

$ws.Range("A24").Value = "L"
$ws.Range("B24").Value = 2

$ws.Range("A25").Value = "r"
$ws.Range("B25").Value = 81

$ws.Range("A26").Value = "v"
$ws.Range("B26").Value = 21

$ws.Range("A27").Value = "p"
$ws.Range("B27").Value = 25

$ws.Range("A28").Value = " "
$ws.Range("B28").Value = 233

$ws.Range("A29").Value = "f"
$ws.Range("B29").Value = 8

$ws.Range("A30").Value = "V"
$ws.Range("B30").Value = 5

$ws.Range("A31").Value = "j"
$ws.Range("B31").Value = 1

$ws.Range("A32").Value = "i"
$ws.Range("B32").Value = 127

$ws.Range("A33").Value = "o"
$ws.Range("B33").Value = 64

$ws.Range("A34").Value = "b"
$ws.Range("B34").Value = 18

$ws.Range("A35").Value = "C"
$ws.Range("B35").Value = 2

$ws.Range("A36").Value = "I"
$ws.Range("B36").Value = 3

$ws.Range("A37").Value = "S"
$ws.Range("B37").Value = 2

$ws.Range("A38").Value = "X"
$ws.Range("B38").Value = 1

$ws.Range("A39").Value = "x"
$ws.Range("B39").Value = 3
